# Update the "Resources" sheet date values in column C (rows 16-27) to 42180
$wb = $excel.ActiveWorkbook
$wsResources = $wb.Worksheets.Item("Resources")
$wsRequirement = $wb.Worksheets.Item("Requirement")

for ($r = 16; $r -le 27; $r++) {
    $wsResources.Cells.Item($r, 3).Value = 42180
}

# Make "Resources" sheet the active sheet/tab, with a new selection
$wsResources.Activate()
$wsResources.Range("C32").Select()

# "Requirement" sheet selection remains G20 (tabSelected flag moves away from it)
$wsRequirement.Range("G20").Select()

# Re-activate Resources sheet so it's the one shown/active when saved
$wsResources.Activate()
